# Generate Report for Handback
# Update the timestamp cells on the Overview, zh-cn, and de-de sheets to
# reflect the latest handoff/handback generation times.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first row (G2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-04 13:09:20"

# zh-cn sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-04 13:09:16"
$wsZhCn.Range("K2").Value = "2016-09-04 13:09:59"

# de-de sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-04 13:09:20"
$wsDeDe.Range("K2").Value = "2016-09-04 13:10:13"
